$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow edits, then re-apply protection.
$ws.Unprotect("lido")

$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.254746924910335
$ws.Range("E2").Value = 0.009883921388346151

$ws.Range("D3").Value = 0.2539283217702668
$ws.Range("E3").Value = 0.005321979776476971

$ws.Range("D4").Value = 0.2446996840093019
$ws.Range("E4").Value = 0.01055155875299763

$ws.Range("D5").Value = 0.2466250693100963
$ws.Range("E5").Value = 0.007594284978761801

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.008324204124784629

$ws.Protect("lido", $true, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false, $false, $false)
